# Clarification on exam extensions:
# Appends a new "Note: ..." explanation about deadline extensions to the
# paragraph that currently ends with "...got this document." (the
# paragraph stating the students have 72 hours to complete the project).

$d = $word.ActiveDocument

# Locate the end of the existing paragraph robustly via Find, instead of a
# hard-coded paragraph index.
$anchor = $d.Content.Duplicate
$found = $anchor.Find.Execute(
    "See the details of submission deadline from where you got this document.",
    $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the anchor paragraph for the extension clarification note."
}

# Collapse to the end of the found text (right before the paragraph mark)
# so subsequent InsertAfter calls append within the same paragraph.
$r = $anchor.Duplicate
$r.Collapse(0)

$r.InsertAfter(" ")
$r.Collapse(0)

$r.InsertAfter("Note: as a general rule, usually there is no deadline extension on this type of exams. And, even if administration grants an extension (e.g., for medical reasons), it should be no more than 50% of the original amount (i.e., a total of ")
$r.Collapse(0)

$r.InsertAfter("108")
$r.Collapse(0)

$r.InsertAfter(" hours in a ")
$r.Collapse(0)

$r.InsertAfter("72")
$r.Collapse(0)

$r.InsertAfter(" hour exam). If for any reason you got granted an extension longer than that, you must contact administration to verify the course responsible had agreed on such extension (there were cases in the past in which such unauthorized extensions were given by mistake). Do ")
$r.Collapse(0)

$notStart = $r.Start
$r.InsertAfter("NOT")
$notEnd = $r.End
$r.Collapse(0)

$r.InsertAfter(" contact the course responsible directly, as exams must be marked anonymously. To make the exam conditions fair to all students, submissions with long extensions that were not authorized by the course responsible will be automatically evaluated as failed (i.e., an ")
$r.Collapse(0)

$fStart = $r.Start
$r.InsertAfter("F")
$fEnd = $r.End
$r.Collapse(0)

$r.InsertAfter(").")
$r.Collapse(0)

# Make "NOT" and the final "F" bold, matching the surrounding emphasis used
# elsewhere in the document for similar warnings.
$boldNot = $d.Range($notStart, $notEnd)
$boldNot.Bold = 1

$boldF = $d.Range($fStart, $fEnd)
$boldF.Bold = 1
